$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A9").Value = "oussama"
$ws.Range("B9").Value = "98821616Oo"
